# Commit: "Fixing name of Sectors to be alligned with Baseline"
#
# The three renewable-energy "Sector" rows (5=wind-offshore row, 6=wind-onshore
# row, 7=PV row) on every year sheet need their labels re-aligned with the
# Baseline workbook's naming, and the numeric data in column E (the only
# non-zero data column for these rows) needs to follow the same row, i.e. the
# data that used to sit under "Onshore wind" now belongs to the row relabelled
# "Onshore wind plants", etc. Net effect on column E is a cyclic shift:
#   new E5 = old E6, new E6 = old E7, new E7 = old E5
#
# This applies uniformly across every year worksheet (2000..2100).

$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count

for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- rotate the column E figures for rows 5/6/7 ---
    $e5 = $ws.Range("E5").Value2
    $e6 = $ws.Range("E6").Value2
    $e7 = $ws.Range("E7").Value2

    $ws.Range("E5").Value2 = $e6
    $ws.Range("E6").Value2 = $e7
    $ws.Range("E7").Value2 = $e5

    # --- relabel the Sector column (C) for rows 5/6/7 ---
    $ws.Range("C5").Value2 = "Onshore wind plants"
    $ws.Range("C6").Value2 = "Photovoltaic plants"
    $ws.Range("C7").Value2 = "Offshore wind plants"
}
